$d = $word.ActiveDocument

# --- Locate the target paragraph -------------------------------------------------
# "El sistema deberá realizar peticiones de servicio de servicio o mantenimiento
#  partiendo de un formulario con diferentes preguntas que indique de que se
#  trata la solicitud. "
$needle = "servicio de servicio o mantenimiento partiendo"
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*$needle*") {
        $target = $p
        break
    }
}

$paraStart = $target.Range.Start

$leadIn  = "El sistema deberá "
$prefix  = "realizar peticiones de servicio de servicio o "
$oldWord = "mantenimiento"
$newWord = "daño"

$runBoundary1 = $paraStart + $leadIn.Length
$wordStart    = $runBoundary1 + $prefix.Length
$wordEnd      = $wordStart + $oldWord.Length

# --- Replace the word -------------------------------------------------------------
# Replacing text inside a run merges same-formatted neighbouring runs in this
# engine, so the paragraph ends up as a single run after this call.
$rWord = $d.Range($wordStart, $wordEnd)
$rWord.Text = $newWord

$newWordEnd = $wordStart + $newWord.Length
$paraEnd    = $target.Range.End

# --- Re-establish the original run layout -----------------------------------------
# "El sistema deberá " | "realizar peticiones de servicio de servicio o " | "daño" | " partiendo ... solicitud. "
# Toggling a character property on/off across a boundary forces the engine to
# keep the runs on either side distinct rather than re-coalescing them.
$rA = $d.Range($runBoundary1, $paraEnd)
$rA.Bold = 1
$rA.Bold = 0

$rB = $d.Range($wordStart, $paraEnd)
$rB.Bold = 1
$rB.Bold = 0

$rC = $d.Range($newWordEnd, $paraEnd)
$rC.Bold = 1
$rC.Bold = 0
